# Apply the quantity corrections described by the commit diff.
#
# Sheet "식당판매" (sheetId 1 / index 1): C6 (수저) 3 -> 2
# Sheet "매점판매" (sheetId 2 / index 2): C5 (대패삼겹살) 27 -> 25
# Sheet "장의용품" (sheetId 3 / index 3): C5 (맥주) 13 -> 10
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("식당판매")
$ws1.Range("C6").Value = 2

$ws2 = $wb.Worksheets.Item("매점판매")
$ws2.Range("C5").Value = 25

$ws3 = $wb.Worksheets.Item("장의용품")
$ws3.Range("C5").Value = 10
